$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 with the values that used to be in row 5 (keep C4/D4/G4/O4 unchanged).
# Leading apostrophe forces these to stay plain text instead of being
# auto-converted to numbers/dates by Excel (J4 looks numeric, K4/L4 look
# like date/time).
$ws.Range("A4").Value = "Name3"
$ws.Range("B4").Value = "Username3"
$ws.Range("E4").Value = "email3@gmail.com"
$ws.Range("F4").Value = "School3"
$ws.Range("H4").Value = "Siem Reap"
$ws.Range("I4").Value = "dom.jpg"
$ws.Range("J4").Value = "'25605"
$ws.Range("K4").Value = "'11/05/2025"
$ws.Range("L4").Value = "'21:22:22"
$ws.Range("M4").Value = "B402"
$ws.Range("N4").Value = "School3"

# Remove the now-duplicate last row (old row 5)
$ws.Rows.Item(5).Delete()
